# Update excess mortality plots - Week 41 (2022)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("excess_mortality_provinces")

# --- Revisions to previously reported "actual deaths" figures (weeks 10-40 of 2022) ---
$ws.Range("U116").Value = 384

$ws.Range("T128").Value = 55

$ws.Range("Z132").Value = 501

$ws.Range("X134").Value = 626

$ws.Range("Z135").Value = 453

$ws.Range("W136").Value = 512
$ws.Range("X136").Value = 600

$ws.Range("T137").Value = 41

$ws.Range("Z138").Value = 448

$ws.Range("U139").Value = 382

$ws.Range("S141").Value = 203

$ws.Range("U142").Value = 372
$ws.Range("W142").Value = 455
$ws.Range("Z142").Value = 451

$ws.Range("X143").Value = 603
$ws.Range("Y143").Value = 61
$ws.Range("Z143").Value = 418

$ws.Range("Q144").Value = 130
$ws.Range("S144").Value = 206
$ws.Range("W144").Value = 483
$ws.Range("X144").Value = 604
$ws.Range("Z144").Value = 426
$ws.Range("AA144").Value = 220

$ws.Range("U145").Value = 436
$ws.Range("W145").Value = 478
$ws.Range("X145").Value = 651
$ws.Range("Y145").Value = 76
$ws.Range("Z145").Value = 430

# --- Week 40 (2022) actuals revised ---
$ws.Range("P146").Value = 128
$ws.Range("Q146").Value = 134
$ws.Range("S146").Value = 241
$ws.Range("T146").Value = 50
$ws.Range("U146").Value = 445
$ws.Range("V146").Value = 215
$ws.Range("W146").Value = 490
$ws.Range("X146").Value = 647
$ws.Range("Y146").Value = 85
$ws.Range("Z146").Value = 498
$ws.Range("AA146").Value = 218

# --- New data: Week 41 (2022) ---
$ws.Range("N147").Value = 2022
$ws.Range("O147").Value = 41
$ws.Range("P147").Value = 101
$ws.Range("Q147").Value = 121
$ws.Range("R147").Value = 107
$ws.Range("S147").Value = 238
$ws.Range("T147").Value = 56
$ws.Range("U147").Value = 464
$ws.Range("V147").Value = 218
$ws.Range("W147").Value = 501
$ws.Range("X147").Value = 672
$ws.Range("Y147").Value = 87
$ws.Range("Z147").Value = 470
$ws.Range("AA147").Value = 247
$ws.Range("AC147").Value = 2022
$ws.Range("AD147").Value = 41

$ws.Range("AE147").Formula = "=ROUND((P147-B147)/B147*100,2)"
$ws.Range("AF147").Formula = "=ROUND((Q147-C147)/C147*100,2)"
$ws.Range("AG147").Formula = "=ROUND((R147-D147)/D147*100,2)"
$ws.Range("AH147").Formula = "=ROUND((S147-E147)/E147*100,2)"
$ws.Range("AI147").Formula = "=ROUND((T147-F147)/F147*100,2)"
$ws.Range("AJ147").Formula = "=ROUND((U147-G147)/G147*100,2)"
$ws.Range("AK147").Formula = "=ROUND((V147-H147)/H147*100,2)"
$ws.Range("AL147").Formula = "=ROUND((W147-I147)/I147*100,2)"
$ws.Range("AM147").Formula = "=ROUND((X147-J147)/J147*100,2)"
$ws.Range("AN147").Formula = "=ROUND((Y147-K147)/K147*100,2)"
$ws.Range("AO147").Formula = "=ROUND((Z147-L147)/L147*100,2)"
$ws.Range("AP147").Formula = "=ROUND((AA147-M147)/M147*100,2)"

# --- Update sheet view / selection to reflect where the user ended up working ---
$ws.Application.ActiveWindow.ScrollRow = 112
$ws.Range("AK146").Select() | Out-Null
